$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F (想去人数 / want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 191
$ws1.Cells.Item(4, 6).Value = 152
$ws1.Cells.Item(5, 6).Value = 1315
$ws1.Cells.Item(6, 6).Value = 18287
$ws1.Cells.Item(7, 6).Value = 373
$ws1.Cells.Item(8, 6).Value = 264
$ws1.Cells.Item(9, 6).Value = 1071
$ws1.Cells.Item(10, 6).Value = 6878
$ws1.Cells.Item(11, 6).Value = 691
$ws1.Cells.Item(12, 6).Value = 161
$ws1.Cells.Item(14, 6).Value = 114
$ws1.Cells.Item(15, 6).Value = 69
$ws1.Cells.Item(16, 6).Value = 219
$ws1.Cells.Item(17, 6).Value = 162
$ws1.Cells.Item(19, 6).Value = 247
$ws1.Cells.Item(20, 6).Value = 58
$ws1.Cells.Item(21, 6).Value = 658
$ws1.Cells.Item(24, 6).Value = 35
$ws1.Cells.Item(26, 6).Value = 994
$ws1.Cells.Item(28, 6).Value = 5172
$ws1.Cells.Item(29, 6).Value = 540
$ws1.Cells.Item(30, 6).Value = 44
$ws1.Cells.Item(32, 6).Value = 75
$ws1.Cells.Item(33, 6).Value = 12116
$ws1.Cells.Item(34, 6).Value = 1288
$ws1.Cells.Item(36, 6).Value = 210
$ws1.Cells.Item(38, 6).Value = 3927

# Sheet "全部类型" (sheet4) - update column F (想去人数 / want-to-go count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 191
$ws4.Cells.Item(4, 6).Value = 152
$ws4.Cells.Item(5, 6).Value = 1315
$ws4.Cells.Item(6, 6).Value = 18288
$ws4.Cells.Item(7, 6).Value = 373
$ws4.Cells.Item(8, 6).Value = 264
$ws4.Cells.Item(9, 6).Value = 1071
$ws4.Cells.Item(10, 6).Value = 6878
$ws4.Cells.Item(12, 6).Value = 161
$ws4.Cells.Item(14, 6).Value = 114
$ws4.Cells.Item(15, 6).Value = 69
$ws4.Cells.Item(16, 6).Value = 219
$ws4.Cells.Item(17, 6).Value = 162
$ws4.Cells.Item(19, 6).Value = 247
$ws4.Cells.Item(20, 6).Value = 58
$ws4.Cells.Item(21, 6).Value = 658
$ws4.Cells.Item(24, 6).Value = 35
$ws4.Cells.Item(26, 6).Value = 994
$ws4.Cells.Item(28, 6).Value = 5172
$ws4.Cells.Item(29, 6).Value = 540
$ws4.Cells.Item(32, 6).Value = 44
$ws4.Cells.Item(34, 6).Value = 75
$ws4.Cells.Item(35, 6).Value = 12116
$ws4.Cells.Item(36, 6).Value = 1288
$ws4.Cells.Item(38, 6).Value = 210
$ws4.Cells.Item(40, 6).Value = 3927
